$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 is a duplicate of row 3, except the id (column A) and CodigoGfh (column J).
$ws.Range("A3:J3").Copy($ws.Range("A6:J6"))
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 10).Value = 1

# Row 5 is a duplicate of row 4, except the id (A), CodiHospi (B), FechaHoraAlta (I) and CodigoGfh (J).
$ws.Range("A4:J4").Copy($ws.Range("A5:J5"))
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "AH582"
$ws.Cells.Item(5, 9).Value = "2020-11-09  20:30:20"
$ws.Cells.Item(5, 10).Value = 2

$excel.CutCopyMode = 0

# Move the selection to the new last row, matching the edited workbook.
$ws.Range("A6").Select()
